$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Red color used to flag "do not change" addresses: RGB(0xC9,0x21,0x1E)
$red = 0x1E21C9

# Rows whose C:H (address/value) columns must not be changed by the user -
# they get highlighted in red, while keeping their existing thin border.
$rowsToFlag = @(12, 14, 16, 17, 19, 20)

foreach ($r in $rowsToFlag) {
    $ws.Range("C" + $r + ":H" + $r).Font.Color = $red
    $ws.Range("I" + $r).Font.Color = $red
}

# New warning banner under the table.
$ws.Range("C26").Value2 = "Do not change values in addresses in red."
$ws.Range("C26").Font.Color = $red
$ws.Range("C26").Font.Size = 20
$ws.Rows.Item(26).RowHeight = 24.45

$ws.Range("D26").ClearContents()

$ws.Range("C27").Value2 = "this may cause the device to become damaged or inoperable. (SEEED mail)"
$ws.Range("C27").Font.Color = $red

$ws.Range("E31").Select()
